$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.777.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.29%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.300.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.46%  "

$ws.Range("E4").Value = "  -0.61%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.53%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +12.54%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.569"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.89%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.46%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.530"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.71%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0804"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.58%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.83%  "

$ws.Range("E13").Value = "  +0.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.651.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.303.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.86%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.84%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.820"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.94%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "46.766.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.43%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +23.12%  "

$ws.Range("E20").Value = "  +4.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.55%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.57%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.88%  "

$ws.Range("E26").Value = "  -1.31%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "44.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +14.49%  "

$ws.Range("E28").Value = "  +1.64%  "

$ws.Range("E29").Value = "  +6.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "147.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.82%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0799"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.21%  "

$ws.Range("E34").Value = "  +3.34%  "

$ws.Range("E35").Value = "  +12.86%  "

$ws.Range("E36").Value = "  +11.33%  "

$ws.Range("E37").Value = "  +2.98%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.97%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.19"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +23.29%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.44%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.62%  "

$ws.Range("E42").Value = "  -0.35%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.62%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.866.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.55%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +11.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "87.96"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +19.24%  "

$ws.Range("E47").Value = "  +10.16%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "74.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +11.75%  "

$ws.Range("E49").Value = "  +9.58%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "97.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.32%  "
